$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 2 (shifts existing rows 2-12 down to 3-13)
$ws.Rows.Item(2).Insert()

# Date-like / numeric-looking text values must be forced to Text format
# before assignment so they are written as shared strings, not date/number
# serials. Reset the style back to Normal afterwards so the row carries no
# explicit style index (matches the plain data rows below it).
$dateCells = "A2","D2","E2"
foreach ($addr in $dateCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("A2").Value = "2024-04-23"
$ws.Range("B2").Value = "민테크"
$ws.Range("C2").Value = "KB"
$ws.Range("D2").Value = "2024-04-26"
$ws.Range("E2").Value = "2024-05-03"
$ws.Range("F2").Value = 31500000
$ws.Range("G2").Value = 3000000
$ws.Range("H2").Value = "-"
$ws.Range("I2").Value = 6500
$ws.Range("J2").Value = 8500
$ws.Range("K2").Value = "-"
$ws.Range("L2").Value = 10500
$ws.Range("M2").Value = "-"
$ws.Range("N2").Value = "-"
$ws.Range("O2").Value = 0
$ws.Range("P2").Value = "-"
$ws.Range("Q2").Value = "-"
$ws.Range("R2").Value = "1529.43 : 1"
$ws.Range("S2").Value = "-"
$ws.Range("T2").Value = "-"

# Strip the inherited style (and the temporary Text number format) from the
# whole new row so it matches the unstyled data rows elsewhere in the sheet.
$ws.Range("A2:T2").Style = "Normal"
